$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each pair of rows had their match-detail columns (B through AB) swapped between them,
# while columns A (row index), C (Div), D (Date) stayed put.
$pairs = @(
    @(73, 74),
    @(114, 115),
    @(122, 123),
    @(159, 160),
    @(173, 174)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AB$r1")
    $range2 = $ws.Range("B$r2`:AB$r2")

    $vals1 = $range1.Value()
    $vals2 = $range2.Value()

    $range1.Value = $vals2
    $range2.Value = $vals1
}
